# live_trading_results.xlsx — Trade #4 closed at 2026-02-17 12:26:42
# Updates Summary + Strategy Status roll-up numbers and appends the new
# trade row to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet roll-up metrics
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.01   # Current Capital
$summary.Range("B4").Value = 0.01      # Total P&L $
$summary.Range("B5").Value = 0.05      # Total P&L %
$summary.Range("B6").Value = 4         # Total Trades
$summary.Range("B7").Value = 2         # Winning Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet — MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.01     # Capital
$status.Range("D4").Value = 4          # Trades
$status.Range("E4").Value = 0.01       # P&L $
$status.Range("F4").Value = 0.01       # P&L %
$status.Range("G4").Value = 50         # Win Rate %

# ---------------------------------------------------------------
# Append the new trade row (#4) to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A5").Value = 4

    # Date/time columns look numeric to Excel's auto-detect, so force
    # them to text first (matches how the existing rows store them as
    # plain strings, not date/time serials), then clear the resulting
    # "stored as text" style back to Normal so no stray format sticks.
    $ws.Range("B5").NumberFormat = "@"
    $ws.Range("B5").Value = "2026-02-17"
    $ws.Range("B5").Style = "Normal"

    $ws.Range("C5").Value = "12:26:36"

    $ws.Range("D5").Value = "MarketMaking"
    $ws.Range("E5").Value = "DOWN"
    $ws.Range("F5").Value = 0.8100000000000001
    $ws.Range("G5").Value = 0.84
    $ws.Range("H5").Value = "CLOSED"
    $ws.Range("I5").Value = 3.7037
    $ws.Range("J5").Value = 0.03
    $ws.Range("K5").Value = 100.01
    $ws.Range("L5").Value = 0
    $ws.Range("M5").Value = 0
    $ws.Range("N5").Value = 0.6
    $ws.Range("O5").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P5").Value = "early_exit"
    $ws.Range("Q5").Value = 0.13
}
